$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.254.88"
$ws.Range("E2").Value = "  -2.01%  "

$ws.Range("D3").Value = "2.189.56"
$ws.Range("E3").Value = "  -6.87%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'296.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.27%  "

$ws.Range("D6").Value = "'81.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.12%  "

$ws.Range("D7").Value = "'0.509"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.00%  "

$ws.Range("E8").Value = "  -0.15%  "

$ws.Range("D9").Value = "'0.466"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.85%  "

$ws.Range("D10").Value = "'0.0768"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.51%  "

$ws.Range("D11").Value = "'29.08"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.27%  "

$ws.Range("D12").Value = "'47.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -10.24%  "

$ws.Range("E13").Value = "  -2.40%  "

$ws.Range("D14").Value = "'6.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.63%  "

$ws.Range("D15").Value = "2.530.84"
$ws.Range("E15").Value = "  -6.88%  "

$ws.Range("D16").Value = "'13.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.17%  "

$ws.Range("D17").Value = "2.178.70"
$ws.Range("E17").Value = "  -7.96%  "

$ws.Range("D18").Value = "'0.708"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.64%  "

$ws.Range("D19").Value = "39.121.26"
$ws.Range("E19").Value = "  -2.56%  "

$ws.Range("D20").Value = "0.0₃0869"
$ws.Range("E20").Value = "  -3.66%  "

$ws.Range("E21").Value = "  -6.54%  "

$ws.Range("D22").Value = "'64.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.51%  "

$ws.Range("D23").Value = "'10.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.08%  "

$ws.Range("D24").Value = "'225.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.60%  "

$ws.Range("E25").Value = "  +0.17%  "

$ws.Range("E26").Value = "  -6.22%  "

$ws.Range("D27").Value = "'1.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("D28").Value = "'22.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.18%  "

$ws.Range("E29").Value = "  -1.92%  "

$ws.Range("D30").Value = "'9.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.80%  "

$ws.Range("D31").Value = "'149.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.83%  "

$ws.Range("D32").Value = "'31.52"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.12%  "

$ws.Range("E33").Value = "  -0.27%  "

$ws.Range("E34").Value = "  -7.07%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0692"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.69%  "

$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'2.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.76%  "

$ws.Range("E37").Value = "  -3.97%  "

$ws.Range("D38").Value = "'0.0961"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.24%  "

$ws.Range("D39").Value = "'15.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.94%  "

$ws.Range("E40").Value = "  -5.70%  "

$ws.Range("D41").Value = "'1.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.22%  "

$ws.Range("E42").Value = "  -5.69%  "

$ws.Range("D43").Value = "1.891.46"
$ws.Range("E43").Value = "  -3.47%  "

$ws.Range("E44").Value = "  -12.14%  "

$ws.Range("D45").Value = "'0.0258"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.78%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'16.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.28%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'8.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.86%  "

$ws.Range("D48").Value = "'2.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.17%  "

$ws.Range("D49").Value = "2.397.19"
$ws.Range("E49").Value = "  -7.58%  "

$ws.Range("D50").Value = "'70.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.57%  "

$ws.Range("D51").Value = "'86.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.17%  "
